$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 carries the workbook default (unstyled) cell style; used below to
# restore each Price cell's style after a temporary Text NumberFormat
# forces the numeric-looking price string (e.g. "6.500", "23.378.05") to
# be stored as a literal string instead of being auto-coerced into a
# Double by COM (which would silently drop significant trailing zeros).
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.378.05'
$ws.Range("D2").Style = $defaultStyle
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.635.23'
$ws.Range("D3").Style = $defaultStyle
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = $defaultStyle
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.27'
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = '  -0.14%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3828'
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").Value = '  +1.32%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '52.08'
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = '  +0.54%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3574'
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = '  -1.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08177'
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = '  +1.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.223'
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = '  -0.64%  '
$ws.Range("E12").Value = '  +0.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.36'
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Value = '  -0.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.423'
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = '  -1.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.302'
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = '  +1.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001232'
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = '  -0.85%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.626.13'
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.06'
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = '  +1.83%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06953'
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = '  +0.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.572'
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = '  +2.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.37'
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = '  -3.04%  '
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.47'
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = '  -1.98%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.375.04'
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.555'
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = '  +4.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.079'
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = '  -4.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.07'
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.85'
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Value = '  +2.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.266'
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = '  -0.53%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.43'
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = '  -0.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.809.35'
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = '  +0.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.084'
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = '  +14.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.148'
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = '  -6.44%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.500'
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = '  -3.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.50'
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = '  +5.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02757'
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = '  -1.95%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2503'
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = '  -0.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.08769'
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = '  -0.48%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.07030'
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = '  -0.83%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.952'
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = '  -2.37%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.346'
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = '  -0.99%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7003'
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = '  -0.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.23'
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.46'
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = '  -4.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6451'
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = '  +0.28%  '
$ws.Range("E46").Value = '  +0.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.277'
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = '  -1.48%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.954'
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = '  -0.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07938'
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = '  -0.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '129.01'
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = '  +2.66%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.187'
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = '  -1.05%  '
